$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.955.59'
$ws.Range("E2").Value = '  +0.21%  '

$ws.Range("D3").Value = '2.650.51'
$ws.Range("E3").Value = '  +1.62%  '

$ws.Range("E4").Value = '  +0.05%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '582.71'
$ws.Range("E5").Value = '  +1.15%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '144.60'
$ws.Range("E6").Value = '  +0.64%  '

$ws.Range("E7").Value = '  +0.16%  '

$ws.Range("E8").Value = '  -0.17%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '6.60'
$ws.Range("E9").Value = '  +0.65%  '

$ws.Range("E10").Value = '  +0.51%  '

$ws.Range("E11").Value = '  +3.12%  '

$ws.Range("E12").Value = '  -0.36%  '

$ws.Range("D13").Value = '3.119.26'
$ws.Range("E13").Value = '  +1.65%  '

$ws.Range("E14").Value = '  +11.13%  '

$ws.Range("D15").Value = '60.918.69'
$ws.Range("E15").Value = '  +0.26%  '

$ws.Range("E16").Value = '  +0.74%  '

$ws.Range("D17").Value = '2.660.65'

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '11.60'
$ws.Range("E18").Value = '  +2.62%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '4.74'
$ws.Range("E19").Value = '  +1.46%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '350.31'
$ws.Range("E20").Value = '  -0.03%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.89'
$ws.Range("E21").Value = '  -0.18%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.997'
$ws.Range("E22").Value = '  -0.07%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.530'
$ws.Range("E23").Value = '  +2.04%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '63.98'
$ws.Range("E24").Value = '  +1.07%  '

$ws.Range("E25").Value = '  +0.78%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.995'
$ws.Range("E26").Value = '  -0.09%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.15'
$ws.Range("E27").Value = '  +3.18%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.00'
$ws.Range("E28").Value = '  +9.19%  '

$ws.Range("E29").Value = '  +1.37%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.85'
$ws.Range("E30").Value = '  +7.67%  '

$ws.Range("E31").Value = '  +0.07%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '164.28'
$ws.Range("E32").Value = '  +1.11%  '

$ws.Range("E33").Value = '  +1.82%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.58'
$ws.Range("E34").Value = '  +7.13%  '

$ws.Range("E35").Value = '  +1.17%  '

$ws.Range("E36").Value = '  +7.37%  '

$ws.Range("B37").Value = 'Stacks'
$ws.Range("C37").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.66'
$ws.Range("E37").Value = '  +3.16%  '

$ws.Range("B38").Value = 'Bittensor'
$ws.Range("C38").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '339.51'
$ws.Range("E38").Value = '  +12.43%  '

$ws.Range("B39").Value = 'SuiNetwork'
$ws.Range("C39").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.917'
$ws.Range("E39").Value = '  +7.81%  '

$ws.Range("B40").Value = 'Filecoin'
$ws.Range("C40").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '4.09'
$ws.Range("E40").Value = '  +4.89%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '38.57'
$ws.Range("E41").Value = '  +1.62%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.26'
$ws.Range("E42").Value = '  +4.26%  '

$ws.Range("B43").Value = 'Mantle'
$ws.Range("C43").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.624'
$ws.Range("E43").Value = '  +2.73%  '

$ws.Range("B44").Value = 'EnergySwap'
$ws.Range("C44").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '20.39'
$ws.Range("E44").Value = '  +1.98%  '

$ws.Range("E45").Value = '  +2.88%  '

$ws.Range("B46").Value = 'Hedera'
$ws.Range("C46").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0564'
$ws.Range("E46").Value = '  +2.41%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '20.66'
$ws.Range("E47").Value = '  +1.67%  '

$ws.Range("B48").Value = 'Aave'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '133.15'
$ws.Range("E48").Value = '  -0.74%  '

$ws.Range("E49").Value = '  +0.99%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.997'
$ws.Range("E50").Value = '  +0.26%  '

$ws.Range("D51").Value = '2.102.95'
$ws.Range("E51").Value = '  +3.72%  '
